# Generate Report for Handoff
#
# The two tracked markdown files swap roles:
#   - cee6e158-e0d4-415f-ad3e-95b85eea2eff.md is now the "in sync" file (was row 2 -> now row 2 data, but originally held by 3e55b075 row)
#   - 3e55b075-fe56-4878-a8e2-bd22e8829534.md now needs a fresh handoff ("Ready for handoff")
#
# Concretely: row 2 on every sheet now describes cee6e158.md, and row 3 now
# describes 3e55b075.md (status/date refreshed + a new "stale handback" error).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"
$overview.Range("B2").Value = "e2e\cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"

$overview.Range("A3").Value = "3e55b075-fe56-4878-a8e2-bd22e8829534.md"
$overview.Range("B3").Value = "e2e\3e55b075-fe56-4878-a8e2-bd22e8829534.md"
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-28 12:49:07"

foreach ($hl in $overview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\3e55b075-fe56-4878-a8e2-bd22e8829534.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"
$zhcn.Range("G2").Value = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.350b4e1ab16e1615031d3b12cf4c507ead14f7a2.zh-cn.xlf"
$zhcn.Range("I2").Value = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"
$zhcn.Range("J2").Value = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.350b4e1ab16e1615031d3b12cf4c507ead14f7a2.zh-cn.xlf"

$zhcn.Range("A3").Value = "3e55b075-fe56-4878-a8e2-bd22e8829534.md"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "3e55b075-fe56-4878-a8e2-bd22e8829534.7c46fd869bf0173b1a5dbb5b11cc31785398ccab.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-28 12:48:59"
$zhcn.Range("I3").Value = "3e55b075-fe56-4878-a8e2-bd22e8829534.md"
$zhcn.Range("J3").Value = "3e55b075-fe56-4878-a8e2-bd22e8829534.7c46fd869bf0173b1a5dbb5b11cc31785398ccab.zh-cn.xlf"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6b53056f7a2a1aa6d8c46b843201e0da660b920/e2e/3e55b075-fe56-4878-a8e2-bd22e8829534.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e6ba56ac3be71d30b63f49e46c9fec7ed395935/e2e/3e55b075-fe56-4878-a8e2-bd22e8829534.md."

foreach ($hl in $zhcn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2' -or $addr -eq '$I$2') {
        $hl.TextToDisplay = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"
    } elseif ($addr -eq '$A$3' -or $addr -eq '$I$3') {
        $hl.TextToDisplay = "3e55b075-fe56-4878-a8e2-bd22e8829534.md"
    }
}

$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"
$dede.Range("G2").Value = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.350b4e1ab16e1615031d3b12cf4c507ead14f7a2.de-de.xlf"
$dede.Range("I2").Value = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"
$dede.Range("J2").Value = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.350b4e1ab16e1615031d3b12cf4c507ead14f7a2.de-de.xlf"

$dede.Range("A3").Value = "3e55b075-fe56-4878-a8e2-bd22e8829534.md"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "3e55b075-fe56-4878-a8e2-bd22e8829534.7c46fd869bf0173b1a5dbb5b11cc31785398ccab.de-de.xlf"
$dede.Range("H3").Value = "2016-08-28 12:49:07"
$dede.Range("I3").Value = "3e55b075-fe56-4878-a8e2-bd22e8829534.md"
$dede.Range("J3").Value = "3e55b075-fe56-4878-a8e2-bd22e8829534.7c46fd869bf0173b1a5dbb5b11cc31785398ccab.de-de.xlf"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6b53056f7a2a1aa6d8c46b843201e0da660b920/e2e/3e55b075-fe56-4878-a8e2-bd22e8829534.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e6ba56ac3be71d30b63f49e46c9fec7ed395935/e2e/3e55b075-fe56-4878-a8e2-bd22e8829534.md."

foreach ($hl in $dede.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2' -or $addr -eq '$I$2') {
        $hl.TextToDisplay = "cee6e158-e0d4-415f-ad3e-95b85eea2eff.md"
    } elseif ($addr -eq '$A$3' -or $addr -eq '$I$3') {
        $hl.TextToDisplay = "3e55b075-fe56-4878-a8e2-bd22e8829534.md"
    }
}

$dede.Columns.Item(16).ColumnWidth = 39.17

Write-Output "Applied handoff-report update"
